$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Práctica Lineal"  (localSheetId = 0)
# ---------------------------------------------------------------------------
$wsLineal = $wb.Worksheets.Item("Práctica Lineal")

$wsLineal.Range("E2").Value = 2.5
$wsLineal.Range("E3").Value = 5

$wsLineal.Range("E6").Formula = "=3*E2+4*E3"

$wsLineal.Range("C9").Formula = "=2*E2+5*E3"
$wsLineal.Range("E9").Value = 30

$wsLineal.Range("C10").Formula = "=4*E2+2*E3"
$wsLineal.Range("E10").Value = 20

$wsLineal.Range("C11").Formula = "=E2"
$wsLineal.Range("E11").Value = 0

$wsLineal.Range("C12").Formula = "=E3"
$wsLineal.Range("E12").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Práctica No Lineal"  (localSheetId = 1)
# ---------------------------------------------------------------------------
$wsNoLineal = $wb.Worksheets.Item("Práctica No Lineal")

$wsNoLineal.Range("E2").Value = 10.466805390419914
$wsNoLineal.Range("E3").Value = -100

$wsNoLineal.Range("E6").Formula = "=(E2^2)+E2*E3"

$wsNoLineal.Range("C9").Formula = "=(E2^3)+E2*E3"
$wsNoLineal.Range("E9").Value = 100

$wsNoLineal.Range("C10").Formula = "=(E2^3)+E2*E3"
$wsNoLineal.Range("E10").Value = 50

$wsNoLineal.Range("C11").Value = -100
$wsNoLineal.Range("E11").Formula = "=E2"

$wsNoLineal.Range("C12").Formula = "=E2"
$wsNoLineal.Range("E12").Value = 100

$wsNoLineal.Range("C13").Value = -100
$wsNoLineal.Range("E13").Formula = "=E3"

$wsNoLineal.Range("C14").Formula = "=E3"
$wsNoLineal.Range("E14").Value = 100

# ---------------------------------------------------------------------------
# Solver parameters stored as hidden workbook-level (sheet-scoped) defined
# names. $wb.Names is ordered exactly like the <definedNames> block in the
# underlying workbook.xml: for each solver_* key, localSheetId=0 (Práctica
# Lineal) comes first, then localSheetId=1 (Práctica No Lineal).
# ---------------------------------------------------------------------------
$wb.Names.Item(1).RefersTo  = "='Práctica Lineal'!`$E`$2:`$E`$3"   # solver_adj (sheet0)
$wb.Names.Item(2).RefersTo  = "='Práctica No Lineal'!`$E`$2:`$E`$3" # solver_adj (sheet1)

$wb.Names.Item(8).RefersTo  = "=1"   # solver_eng (sheet1): 2 -> 1

$wb.Names.Item(24).RefersTo = "=2"   # solver_lin (sheet1): 1 -> 2

$wb.Names.Item(34).RefersTo = "=2"   # solver_neg (sheet1): 1 -> 2

$wb.Names.Item(37).RefersTo = "=4"   # solver_num (sheet0): 0 -> 4
$wb.Names.Item(38).RefersTo = "=6"   # solver_num (sheet1): 0 -> 6

$wb.Names.Item(39).RefersTo = "='Práctica Lineal'!`$E`$6"    # solver_opt (sheet0)
$wb.Names.Item(40).RefersTo = "='Práctica No Lineal'!`$E`$6" # solver_opt (sheet1)

$wb.Names.Item(84).RefersTo = "=2"   # solver_typ (sheet1): 1 -> 2

$wb.Save()
